$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Asset_Cal_Info")

# ---------------------------------------------------------------------------
# NUTNR Cal Coefficients for GPG 514, 528, 566
#   H6  -> CC_cal_temp   (numeric)
#   H7  -> CC_wl         (array, stored as text)
#   H8  -> CC_eno3       (array, stored as text)
#   H9  -> CC_eswa       (array, stored as text)
#   H10 -> CC_di         (array, stored as text)
#   H11 -> CC_lower_wavelength_limit_for_spectra_fit (numeric, new highlighted style)
#   H12 -> CC_upper_wavelength_limit_for_spectra_fit (numeric, new highlighted style)
# ---------------------------------------------------------------------------

$ccWl = @'
[189.69, 190.47, 191.26, 192.04, 192.83, 193.61, 194.39, 195.18, 195.97, 196.75, 197.54, 198.32, 199.11, 199.9, 200.69, 201.47, 202.26, 203.05, 203.84, 204.63, 205.42, 206.2, 206.99, 207.78, 208.57, 209.36, 210.15, 210.95, 211.74, 212.53, 213.32, 214.11, 214.9, 215.7, 216.49, 217.28, 218.07, 218.87, 219.66, 220.45, 221.25, 222.04, 222.84, 223.63, 224.43, 225.22, 226.02, 226.81, 227.61, 228.4, 229.2, 230, 230.79, 231.59, 232.39, 233.18, 233.98, 234.78, 235.58, 236.37, 237.17, 237.97, 238.77, 239.57, 240.37, 241.17, 241.97, 242.77, 243.57, 244.37, 245.17, 245.97, 246.77, 247.57, 248.37, 249.17, 249.97, 250.77, 251.57, 252.37, 253.18, 253.98, 254.78, 255.58, 256.38, 257.19, 257.99, 258.79, 259.59, 260.4, 261.2, 262, 262.81, 263.61, 264.41, 265.22, 266.02, 266.83, 267.63, 268.44, 269.24, 270.04, 270.85, 271.65, 272.46, 273.26, 274.07, 274.87, 275.68, 276.49, 277.29, 278.1, 278.9, 279.71, 280.52, 281.32, 282.13, 282.93, 283.74, 284.55, 285.35, 286.16, 286.97, 287.77, 288.58, 289.39, 290.19, 291, 291.81, 292.62, 293.42, 294.23, 295.04, 295.85, 296.65, 297.46, 298.27, 299.08, 299.88, 300.69, 301.5, 302.31, 303.12, 303.92, 304.73, 305.54, 306.35, 307.16, 307.96, 308.77, 309.58, 310.39, 311.2, 312.01, 312.81, 313.62, 314.43, 315.24, 316.05, 316.86, 317.67, 318.47, 319.28, 320.09, 320.9, 321.71, 322.52, 323.33, 324.13, 324.94, 325.75, 326.56, 327.37, 328.18, 328.99, 329.79, 330.6, 331.41, 332.22, 333.03, 333.84, 334.64, 335.45, 336.26, 337.07, 337.88, 338.69, 339.49, 340.3, 341.11, 341.92, 342.73, 343.53, 344.34, 345.15, 345.96, 346.77, 347.57, 348.38, 349.19, 350, 350.8, 351.61, 352.42, 353.23, 354.03, 354.84, 355.65, 356.45, 357.26, 358.07, 358.87, 359.68, 360.49, 361.29, 362.1, 362.91, 363.71, 364.52, 365.32, 366.13, 366.94, 367.74, 368.55, 369.35, 370.16, 370.96, 371.77, 372.57, 373.38, 374.18, 374.99, 375.79, 376.6, 377.4, 378.21, 379.01, 379.81, 380.62, 381.42, 382.23, 383.03, 383.83, 384.64, 385.44, 386.24, 387.04, 387.85, 388.65, 389.45, 390.25, 391.06, 391.86, 392.66, 393.46, 394.26]
'@

$ccEno3 = @'
[-0.00021425, -0.00461322, -0.03471286, -0.0012909, -0.01002091, 0.01572297, 0.0137133, 0.00178863, -0.00303613, -0.00114568, -0.00196477, -0.00884724, -0.0080948, 0.00074724, -0.01742237, 0.00992815, -0.00428105, 0.01929635, -0.00725399, -0.00201656, -0.00014736, 0.00232121, 0.01006033, 0.00643182, 0.00592234, 0.00827105, 0.00793879, 0.0071555, 0.00681379, 0.00660685, 0.00613918, 0.00588971, 0.00556329, 0.00524898, 0.00487834, 0.00456418, 0.00421566, 0.00389696, 0.00358277, 0.00327313, 0.00297306, 0.00268393, 0.00241168, 0.00214562, 0.00190938, 0.00168818, 0.00146822, 0.00127449, 0.00111369, 0.00095278, 0.00082719, 0.0007016, 0.00056811, 0.00048025, 0.00040517, 0.00033889, 0.00028127, 0.0002294, 0.0001781, 0.00014223, 0.0001068, 0.00008257, 0.00006559, 0.00004768, 0.0000257, 0.0000245, 0.00000502, 0.00000903, -0.00000713, -0.00000045, 0.00000078, -0.00000975, -0.000003, -0.00001411, -0.00001184, -0.00001323, -0.00002696, -0.0000161, -0.00000901, -0.00002004, -0.00000137, 0.00000523, 0.00000359, 0.00000762, 0.00000634, 0.00002188, 0.00000603, 0.00000405, 0.00001322, 0.00001665, 0.00001321, 0.00001993, 0.00001924, 0.00003445, 0.00003655, 0.00003616, 0.00004346, 0.00004745, 0.00005631, 0.0000323, 0.00005591, 0.00003432, 0.00004427, 0.00005948, 0.00005379, 0.00007259, 0.00006885, 0.00006693, 0.0000771, 0.00007856, 0.00008285, 0.00008123, 0.00009142, 0.0000802, 0.00007889, 0.0000932, 0.00009667, 0.0001042, 0.00007186, 0.00009852, 0.00010934, 0.00010417, 0.00009298, 0.00010215, 0.00011711, 0.00013426, 0.00011756, 0.00012793, 0.00014444, 0.00013154, 0.00014673, 0.00013831, 0.00015824, 0.00013712, 0.00014341, 0.00016351, 0.0001582, 0.00015813, 0.00015783, 0.00013898, 0.00015188, 0.0001676, 0.0001659, 0.00016325, 0.00016915, 0.00016691, 0.00018059, 0.00019081, 0.00019126, 0.00018127, 0.00020239, 0.0001831, 0.00019906, 0.00020395, 0.0001953, 0.00020427, 0.00020199, 0.0002171, 0.00020962, 0.00023148, 0.0002201, 0.00021582, 0.0002339, 0.00023631, 0.00022546, 0.00023706, 0.00023951, 0.00026088, 0.00023879, 0.00024515, 0.00023292, 0.00026366, 0.00024974, 0.00025303, 0.00025803, 0.00027393, 0.00027042, 0.00027642, 0.000282, 0.00026678, 0.00027177, 0.00027491, 0.00025061, 0.00027081, 0.00029201, 0.0002908, 0.00028492, 0.0003058, 0.00031599, 0.00031152, 0.00031588, 0.00033289, 0.00028729, 0.00030352, 0.00033357, 0.00033272, 0.00032157, 0.00031632, 0.00033055, 0.00033786, 0.00034263, 0.00034273, 0.00034455, 0.00036475, 0.00035646, 0.00035561, 0.00036635, 0.00035815, 0.0003601, 0.00036304, 0.00037052, 0.0003858, 0.0003736, 0.00037849, 0.00040195, 0.00040337, 0.00039936, 0.00038445, 0.00041048, 0.00042139, 0.00040596, 0.00040546, 0.00041635, 0.00042452, 0.00041937, 0.00043566, 0.00043013, 0.00044192, 0.00045126, 0.00045019, 0.00047114, 0.00043082, 0.0004531, 0.00044963, 0.00048421, 0.00043218, 0.00050073, 0.00048945, 0.00050841, 0.00047831, 0.00050603, 0.00046158, 0.00051902, 0.00049463, 0.00050197, 0.00050619, 0.00053299, 0.00052907, 0.00050227, 0.00053418, 0.00050914, 0.00053763, 0.0005046, 0.00055875, 0.00054377, 0.0005558]
'@

$ccEswa = @'
[-0.00035427, 0.01204174, 0.04747922, 0.00498282, 0.00304351, 0.00153688, 0.0072225, -0.00742801, -0.00881144, -0.00395838, 0.00675199, 0.02111784, 0.03217874, 0.04164911, 0.07202859, 0.06473948, 0.0772996, 0.08350089, 0.08567336, 0.08473379, 0.08024467, 0.07674616, 0.07270172, 0.06528003, 0.05738647, 0.04886592, 0.04034703, 0.03416583, 0.02805719, 0.02298243, 0.0187789, 0.01517481, 0.01221288, 0.00976399, 0.00782712, 0.0062149, 0.00491848, 0.00387425, 0.00304478, 0.00239784, 0.00187893, 0.00146534, 0.00114971, 0.00089933, 0.00071807, 0.00053966, 0.00043759, 0.00033895, 0.00027803, 0.00022638, 0.00015934, 0.0001462, 0.00013217, 0.00011779, 0.00009441, 0.00008257, 0.00006661, 0.00004581, 0.00004529, 0.0000359, 0.00003726, 0.00002604, 0.00002465, 0.00002098, 0.00001952, 0.00000809, 0.00000959, -0.00000734, 0.0000028, 0.0000038, -0.00001327, -0.00000684, -0.00001424, -0.00000525, -0.00000275, 0.0000007, 0.00000327, 0.000015, -0.00000725, -0.00000302, -0.00001884, -0.00000945, -0.00000364, 0.00001001, -0.00001564, -0.00001059, 0.00001512, 0.0000003, 0.00002993, 0.00001409, 0.00002444, 0.00001761, 0.00001659, 0.0000085, 0.00001253, 0.00003227, 0.00002514, 0.00000308, 0.00002011, 0.0000263, 0.00002848, 0.00004113, 0.00003526, 0.00001865, 0.00004321, 0.00002842, 0.00002663, 0.00003159, 0.00004389, 0.0000231, 0.00004319, 0.00001916, 0.00003243, 0.00005018, 0.00003757, 0.0000417, 0.00004012, 0.00002913, 0.00005817, 0.00003749, 0.00005662, 0.00006599, 0.00007015, 0.00005121, 0.00007231, 0.00003452, 0.00004603, 0.00006855, 0.00004657, 0.00005491, 0.00004079, 0.00007177, 0.0000615, 0.00007173, 0.00005326, 0.00005364, 0.00006645, 0.00006987, 0.0000659, 0.00009308, 0.00009397, 0.00009191, 0.00007844, 0.00010098, 0.00009358, 0.00009028, 0.00009908, 0.00008954, 0.00007818, 0.00010396, 0.00008845, 0.00010103, 0.00010066, 0.00010377, 0.0001249, 0.00011409, 0.00011746, 0.00010487, 0.00012044, 0.00012053, 0.00011097, 0.0001315, 0.00010097, 0.00011714, 0.0001158, 0.00014344, 0.00012429, 0.00011603, 0.00014677, 0.00012552, 0.00015669, 0.00012656, 0.00016034, 0.00014542, 0.00015172, 0.00013399, 0.00013722, 0.00013, 0.00014002, 0.00017287, 0.00019514, 0.00017033, 0.0002091, 0.00017748, 0.00017451, 0.00014968, 0.00016483, 0.00017048, 0.00017453, 0.00019243, 0.00021182, 0.00019687, 0.00021466, 0.00020756, 0.00020483, 0.00016977, 0.00019641, 0.00020883, 0.0002028, 0.00019185, 0.00020231, 0.00020758, 0.000219, 0.00022072, 0.00021994, 0.00024393, 0.00022083, 0.00024006, 0.00022563, 0.00023015, 0.00023681, 0.00022133, 0.00024952, 0.00021158, 0.00022606, 0.00021206, 0.00022996, 0.00025155, 0.00025662, 0.00024803, 0.00026479, 0.00025356, 0.0002693, 0.00028489, 0.00026545, 0.0002224, 0.00027264, 0.00024901, 0.00027484, 0.00024837, 0.00029573, 0.00025025, 0.0002356, 0.00025622, 0.00027849, 0.000271, 0.00024874, 0.00026675, 0.00026861, 0.00029677, 0.00027193, 0.00029015, 0.00028966, 0.00030541, 0.00031071, 0.00028235, 0.00027932, 0.00029573, 0.00029888, 0.00035826, 0.00031793, 0.00030131, 0.0003243, 0.0002905, 0.00027398, 0.0002634]
'@

$ccDi = @'
[35.5, 30.25, 32.45833333, 31.375, 37.16666667, 33.16666667, 30.79166667, 35.79166667, 28.625, 28.79166667, 37, 53.125, 178.5, 614.70833333, 1794.58333333, 4187.20833333, 7846.95833333, 12101.54166667, 15879.25, 18519.95833333, 19933.625, 20472.75, 20511.41666667, 20416.75, 20357.875, 20476.625, 20852.29166667, 21491.45833333, 22448.08333333, 23705.375, 25275, 27171.04166667, 29330, 31683.54166667, 34107.125, 36502.91666667, 38640.25, 40346.45833333, 41497.5, 41946.45833333, 41677.625, 40784.04166667, 39421.66666667, 37795.33333333, 36123.79166667, 34484.625, 33101.45833333, 31984, 31139.625, 30600.95833333, 30413.75, 30464.91666667, 30822.58333333, 31463.70833333, 32350.41666667, 33464, 34839.125, 36389.58333333, 38099.95833333, 39976.08333333, 41872.20833333, 43793.83333333, 45632.91666667, 47223.16666667, 48541.79166667, 49454.83333333, 49834.45833333, 49728.54166667, 49063.45833333, 47855.54166667, 46289.125, 44420.20833333, 42358.5, 40296.70833333, 38243.33333333, 36287.08333333, 34509.20833333, 32913.16666667, 31490.20833333, 30271.25, 29225.08333333, 28354.95833333, 27641.33333333, 27090.79166667, 26687.33333333, 26424.83333333, 26312.16666667, 26343.16666667, 26488.375, 26776.79166667, 27182.875, 27730.16666667, 28352.79166667, 29089.625, 29891.54166667, 30729.79166667, 31600.20833333, 32447.08333333, 33196.66666667, 33847.54166667, 34313.83333333, 34538.33333333, 34498.625, 34226, 33656.33333333, 32870.70833333, 31899.41666667, 30787.125, 29602.54166667, 28403.08333333, 27239.41666667, 26124.58333333, 25064.875, 24131.20833333, 23292.04166667, 22561.16666667, 21917.66666667, 21388.70833333, 20958.875, 20624.91666667, 20386.375, 20229.04166667, 20176.08333333, 20203.33333333, 20319.54166667, 20523.875, 20811.16666667, 21183.79166667, 21647.66666667, 22158.66666667, 22759.16666667, 23433.20833333, 24140.33333333, 24892.25, 25708.04166667, 26547.66666667, 27412.625, 28299.125, 29182.125, 30028.875, 30849.70833333, 31582.83333333, 32280.79166667, 32886.66666667, 33348.45833333, 33708.33333333, 33950.875, 34063.375, 34033.5, 33882.125, 33552.41666667, 33145.66666667, 32638.70833333, 32021.16666667, 31353.16666667, 30626.79166667, 29856.16666667, 29071.41666667, 28295.70833333, 27496, 26763.08333333, 26053.79166667, 25376.20833333, 24756.08333333, 24175.83333333, 23643.625, 23155.625, 22695.79166667, 22268.16666667, 21843.875, 21444.95833333, 21048.75, 20665.66666667, 20317.83333333, 20006, 19731.125, 19495.83333333, 19293.41666667, 19140.16666667, 19013.875, 18922.58333333, 18863.20833333, 18818.79166667, 18808.79166667, 18786.54166667, 18760.66666667, 18734.125, 18737.41666667, 18721.375, 18728.16666667, 18718.83333333, 18675.70833333, 18669.83333333, 18681.08333333, 18698.75, 18745.625, 18783.66666667, 18836.29166667, 18873.66666667, 18901.875, 18959, 19022.20833333, 19080.5, 19111, 19149.625, 19152.58333333, 19133.29166667, 19120.79166667, 19070, 19001.95833333, 18897.66666667, 18755.33333333, 18589.375, 18411.375, 18179.33333333, 17938.20833333, 17676.875, 17376.25, 17075.70833333, 16781.45833333, 16496.41666667, 16227.58333333, 15928.75, 15567.66666667, 15175.66666667, 14760.5, 15745, 15352, 15026, 14692, 14489, 14195, 13974, 13752, 13602, 13495, 13454, 13349, 13107, 12834, 12651, 12463, 12330, 12205, 12059, 11928, 11881, 11890, 11806, 11524, 10981, 10364, 9625, 8729, 7637, 6328]
'@

# CC_cal_temp: simple numeric constant
$ws.Range("H6").Value = 19.93

# CC_wl / CC_eno3 / CC_eswa / CC_di: long numeric arrays stored as text
$ws.Range("H7").Value = $ccWl
$ws.Range("H8").Value = $ccEno3
$ws.Range("H9").Value = $ccEswa
$ws.Range("H10").Value = $ccDi

# CC_lower_wavelength_limit_for_spectra_fit / CC_upper_wavelength_limit_for_spectra_fit
# These two get a new highlighted (green fill) numeric style, matching the
# look-and-feel of the neighboring "Sensor OOIBARCODE" column (I11/I12),
# plus a 2-decimal number format. Build the style once on H11, then reuse it
# for H12 by copying the already-styled H11 cell (so only a single new cell
# style entry gets created, instead of one per cell).
$ws.Range("I11").Copy() | Out-Null
$ws.Range("H11").PasteSpecial(-4122) | Out-Null
$ws.Range("H11").Value = 217
$ws.Range("H11").NumberFormat = "0.00"
$ws.Range("H11").Interior.Color = 5296274

$ws.Range("H11").Copy() | Out-Null
$ws.Range("H12").PasteSpecial(-4122) | Out-Null
$ws.Range("H12").Value = 240

# ---------------------------------------------------------------------------
# Sheet/selection activity: Moorings becomes the active tab, and the
# Asset_Cal_Info selection moves to H10.
# ---------------------------------------------------------------------------
$ws.Range("H10").Select() | Out-Null

$wsMoorings = $wb.Worksheets.Item("Moorings")
$wsMoorings.Activate()
